$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.088.53"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.651.18"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5198"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06331"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07678"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.594"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D13").Value = "1.646.62"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "1.878.40"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5588"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "0.0₅8136"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "26.094.35"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.623"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.915"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.223"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05490"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.447"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.360"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9481"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.416"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.781"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5635"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01577"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.853"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "1.029.02"
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "1.791.89"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4336"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.973"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("E51").Value = "  -2.53%  "
